# Refresh the cryptocurrency price / 1h-volume table with the latest scraped values
# (mirrors the periodic "Updated cryptos list ... with GitHub Actions" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.413.92'
$ws.Range('E2').Value = '  -2.92%  '
$ws.Range('D3').Value = '2.452.25'
$ws.Range('E3').Value = '  -3.94%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'529.58"
$ws.Range('E5').Value = '  -1.92%  '
$ws.Range('D6').Value = "'133.99"
$ws.Range('E6').Value = '  -7.12%  '
$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = "'0.556"
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('D9').Value = '2.459.63'
$ws.Range('E9').Value = '  -4.51%  '
$ws.Range('D10').Value = "'0.0990"
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').Value = "'5.30"
$ws.Range('E12').Value = '  -3.35%  '
$ws.Range('E13').Value = '  -5.38%  '
$ws.Range('D14').Value = '2.887.75'
$ws.Range('D15').Value = '58.327.41'
$ws.Range('E15').Value = '  -3.01%  '
$ws.Range('D16').Value = "'22.66"
$ws.Range('E16').Value = '  -5.67%  '
$ws.Range('E17').Value = '  -3.72%  '
$ws.Range('D18').Value = '2.460.43'
$ws.Range('E18').Value = '  -3.81%  '
$ws.Range('D19').Value = "'10.77"
$ws.Range('E19').Value = '  -4.25%  '
$ws.Range('D20').Value = "'4.20"
$ws.Range('E20').Value = '  -3.02%  '
$ws.Range('D21').Value = "'320.90"
$ws.Range('E21').Value = '  -1.85%  '
$ws.Range('D22').Value = "'0.997"
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('E23').Value = '  -4.18%  '
$ws.Range('D24').Value = "'62.52"
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('D25').Value = "'0.408"
$ws.Range('E25').Value = '  -5.83%  '
$ws.Range('E26').Value = '  -2.36%  '
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').Value = "'7.45"
$ws.Range('E28').Value = '  -7.22%  '
$ws.Range('D29').Value = '0.0₃0751'
$ws.Range('E29').Value = '  -5.51%  '
$ws.Range('D30').Value = "'6.51"
$ws.Range('E30').Value = '  -7.74%  '
$ws.Range('D31').Value = "'1.75"
$ws.Range('E31').Value = '  -3.70%  '
$ws.Range('D32').Value = "'164.51"
$ws.Range('E32').Value = '  -1.06%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('E34').Value = '  -6.56%  '
$ws.Range('D35').Value = "'18.23"
$ws.Range('E35').Value = '  -2.70%  '
$ws.Range('E36').Value = '  -8.58%  '
$ws.Range('D37').Value = "'4.02"
$ws.Range('E37').Value = '  -8.44%  '
$ws.Range('E38').Value = '  -6.20%  '
$ws.Range('D39').Value = "'36.45"
$ws.Range('E39').Value = '  -1.84%  '
$ws.Range('E40').Value = '  -3.62%  '
$ws.Range('D41').Value = "'3.55"
$ws.Range('E41').Value = '  -4.78%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = "'274.20"
$ws.Range('E42').Value = '  -9.06%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'0.997"
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').Value = "'5.07"
$ws.Range('E44').Value = '  -9.44%  '
$ws.Range('D45').Value = "'10.82"
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').Value = "'0.585"
$ws.Range('E46').Value = '  -4.40%  '
$ws.Range('D47').Value = "'0.0922"
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').Value = "'121.07"
$ws.Range('E48').Value = '  -5.24%  '
$ws.Range('D49').Value = "'0.0504"
$ws.Range('E49').Value = '  -2.91%  '
$ws.Range('E50').Value = '  -5.27%  '
$ws.Range('D51').Value = "'17.14"
$ws.Range('E51').Value = '  -6.01%  '

Write-Output "Updated 90 cells across the cryptos worksheet"
